$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 186 / 187: the two fixtures' data got swapped (Monterrey
#    U23 vs Mazatlan FC U23 <-> Unam Pumas U23 vs Tijuana U23), while
#    the match-number column (A) and the shared date (E) stay as-is.
# ------------------------------------------------------------------
$swapCols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

foreach ($col in $swapCols) {
  $top = $ws.Cells.Item(186, $col).Value()
  $bot = $ws.Cells.Item(187, $col).Value()
  $ws.Cells.Item(186, $col).Value = $bot
  $ws.Cells.Item(187, $col).Value = $top
}

# ------------------------------------------------------------------
# 2) Two brand-new fixtures are inserted right before the last row
#    (old row 229), pushing it down to row 231 and bumping the
#    dimension from AC229 to AC231.
# ------------------------------------------------------------------
$ws.Rows(229).Insert()
$ws.Rows(229).Insert()

# carry the existing formatting (bold/border on A, date format on E)
# down onto the two freshly inserted rows
$ws.Range("A228").Copy()
$ws.Range("A229:A230").PasteSpecial(-4122)
$ws.Range("E228").Copy()
$ws.Range("E229:E230").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 229: Atlante vs Correcaminos
$ws.Cells.Item(229,1).Value = 227
$ws.Cells.Item(229,2).Value = 7640653
$ws.Cells.Item(229,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(229,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(229,5).Value = 45393.92013888889
$ws.Cells.Item(229,6).Value = "Atlante"
$ws.Cells.Item(229,7).Value = "Correcaminos"
$ws.Cells.Item(229,8).Value = 4
$ws.Cells.Item(229,9).Value = 0
$ws.Cells.Item(229,10).Value = "H"
$ws.Cells.Item(229,11).Value = 1.6
$ws.Cells.Item(229,12).Value = 3.8
$ws.Cells.Item(229,13).Value = 4.5
$ws.Cells.Item(229,14).Value = 1.4
$ws.Cells.Item(229,15).Value = 4.5
$ws.Cells.Item(229,16).Value = 8
$ws.Cells.Item(229,17).Value = -1.25
$ws.Cells.Item(229,18).Value = 1.85
$ws.Cells.Item(229,19).Value = 1.95
$ws.Cells.Item(229,20).Value = 2.5
$ws.Cells.Item(229,21).Value = 1.8
$ws.Cells.Item(229,22).Value = 2
$ws.Cells.Item(229,23).Value = 0.3999999999999999
$ws.Cells.Item(229,24).Value = -1
$ws.Cells.Item(229,25).Value = -1
$ws.Cells.Item(229,26).Value = 0.8500000000000001
$ws.Cells.Item(229,27).Value = -1
$ws.Cells.Item(229,28).Value = 0.8
$ws.Cells.Item(229,29).Value = -1

# New row 230: Club Celaya vs Tapatio
$ws.Cells.Item(230,1).Value = 228
$ws.Cells.Item(230,2).Value = 7641727
$ws.Cells.Item(230,3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(230,4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(230,5).Value = 45394.00347222222
$ws.Cells.Item(230,6).Value = "Club Celaya"
$ws.Cells.Item(230,7).Value = "Tapatio"
$ws.Cells.Item(230,8).Value = 1
$ws.Cells.Item(230,9).Value = 2
$ws.Cells.Item(230,10).Value = "A"
$ws.Cells.Item(230,11).Value = 1.8
$ws.Cells.Item(230,12).Value = 3.4
$ws.Cells.Item(230,13).Value = 3.8
$ws.Cells.Item(230,14).Value = 2.375
$ws.Cells.Item(230,15).Value = 3.6
$ws.Cells.Item(230,16).Value = 2.875
$ws.Cells.Item(230,17).Value = -0.25
$ws.Cells.Item(230,18).Value = 2.025
$ws.Cells.Item(230,19).Value = 1.775
$ws.Cells.Item(230,20).Value = 2.75
$ws.Cells.Item(230,21).Value = 1.925
$ws.Cells.Item(230,22).Value = 1.775
$ws.Cells.Item(230,23).Value = -1
$ws.Cells.Item(230,24).Value = -1
$ws.Cells.Item(230,25).Value = 1.875
$ws.Cells.Item(230,26).Value = -1
$ws.Cells.Item(230,27).Value = 0.7749999999999999
$ws.Cells.Item(230,28).Value = 0.4625
$ws.Cells.Item(230,29).Value = -0.5

# ------------------------------------------------------------------
# 3) Row 231 (the former row 229, Oaxaca vs Atletico Morelia) keeps
#    its match id (column B) but the running counter in column A
#    advances to match its new row position, and several of its odds
#    columns were refreshed.
# ------------------------------------------------------------------
$ws.Cells.Item(231,1).Value = 229
$ws.Cells.Item(231,14).Value = 2.375
$ws.Cells.Item(231,15).Value = 3.6
$ws.Cells.Item(231,17).Value = 0
$ws.Cells.Item(231,18).Value = 1.775
$ws.Cells.Item(231,19).Value = 2.025
